$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, matching the style of H1 (bold header with border)
$ws.Cells.Item(1, 9).Value2 = "I0"
$ws.Cells.Item(1, 10).Value2 = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate I and J columns (I0 and IF) for data rows 2-86
$data = @{
    2 = @(8, 8)
    3 = @(9, 9)
    4 = @(9, 9)
    5 = @(9, 9)
    6 = @(8, 8)
    7 = @(8, 8)
    8 = @(9, 9)
    9 = @(9, 9)
    10 = @(8, 9)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(9, 9)
    15 = @(9, 9)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(9, 9)
    20 = @(9, 9)
    21 = @(9, 9)
    22 = @(8, 8)
    23 = @(9, 9)
    24 = @(9, 9)
    25 = @(8, 8)
    26 = @(9, 9)
    27 = @(9, 9)
    28 = @(9, 9)
    29 = @(9, 9)
    30 = @(9, 9)
    31 = @(9, 9)
    32 = @(9, 9)
    33 = @(9, 10)
    34 = @(8, 8)
    35 = @(9, 9)
    36 = @(9, 9)
    37 = @(9, 9)
    38 = @(8, 8)
    39 = @(9, 9)
    40 = @(9, 9)
    41 = @(9, 9)
    42 = @(9, 9)
    43 = @(9, 9)
    44 = @(8, 8)
    45 = @(9, 9)
    46 = @(9, 9)
    47 = @(9, 9)
    48 = @(9, 9)
    49 = @(9, 9)
    50 = @(9, 9)
    51 = @(9, 9)
    52 = @(9, 9)
    53 = @(9, 9)
    54 = @(8, 9)
    55 = @(9, 9)
    56 = @(9, 9)
    57 = @(9, 9)
    58 = @(10, 10)
    59 = @(9, 9)
    60 = @(9, 9)
    61 = @(9, 9)
    62 = @(9, 9)
    63 = @(8, 9)
    64 = @(9, 9)
    65 = @(9, 9)
    66 = @(8, 8)
    67 = @(10, 10)
    68 = @(7, 7)
    69 = @(8, 8)
    70 = @(10, 10)
    71 = @(8, 9)
    72 = @(9, 9)
    73 = @(8, 9)
    74 = @(9, 9)
    75 = @(9, 9)
    76 = @(9, 9)
    77 = @(9, 9)
    78 = @(8, 9)
    79 = @(7, 7)
    80 = @(6, 6)
    81 = @(7, 7)
    82 = @(5, 5)
    83 = @(7, 7)
    84 = @(4, 4)
    85 = @(4, 4)
    86 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item([int]$row, 9).Value2 = $vals[0]
    $ws.Cells.Item([int]$row, 10).Value2 = $vals[1]
}
